# Update numeric "want to go" counts (column F) on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 6650
$wsExhibit.Range("F8").Value = 135
$wsExhibit.Range("F9").Value = 6101
$wsExhibit.Range("F12").Value = 1244
$wsExhibit.Range("F13").Value = 1244
$wsExhibit.Range("F14").Value = 9
$wsExhibit.Range("F15").Value = 91
$wsExhibit.Range("F19").Value = 357
$wsExhibit.Range("F22").Value = 4411
$wsExhibit.Range("F26").Value = 39

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 6650
$wsAll.Range("F8").Value = 135
$wsAll.Range("F9").Value = 6101
$wsAll.Range("F12").Value = 1244
$wsAll.Range("F13").Value = 1244
$wsAll.Range("F14").Value = 9
$wsAll.Range("F15").Value = 91
$wsAll.Range("F19").Value = 357
$wsAll.Range("F22").Value = 4411
$wsAll.Range("F27").Value = 39
